$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 16.4.2.1: collapse the 4-category table (Voluntarily surrendered /
# Seized / Abducted / Lost) down to just 2 categories (Voluntarily
# surrendered / Seized), and extend the year range with a new 2020
# column.
# ------------------------------------------------------------------

# 1) Drop the "Abducted" (Похищенные/Уурдалган) and "Lost"
#    (Утерянные/Жоготулган) rows entirely - this shifts nothing below
#    them since they were the last two data rows.
$ws.Rows("6:7").Delete()

# 2) Row 5 ("Seized firearms" row) is now the last row of the table,
#    so it needs to pick up the heavier "closing" bottom-border format
#    that row 7 used to have (now vacated). Re-use row 2's existing
#    format (same bottom border) instead of fabricating a new style.
$ws.Range("A2:G2").Copy()
$ws.Range("A5:G5").PasteSpecial(-4122) | Out-Null

# 3) Add the new year column (H) mirroring the existing D:G columns.
#    Copy formats from column G so the new column reuses the same
#    underlying cell styles rather than creating new ones.
$ws.Range("G2").Copy()
$ws.Range("H2").PasteSpecial(-4122) | Out-Null

$ws.Range("G3").Copy()
$ws.Range("H3").PasteSpecial(-4122) | Out-Null

$ws.Range("G4").Copy()
$ws.Range("H4").PasteSpecial(-4122) | Out-Null

$ws.Range("H5").Copy()
# (H5 already picked up row 5's new format from step 2 via A2:G2->A5:G5
#  copy not covering H5, so explicitly mirror the closing-row format.)
$ws.Range("G5").Copy()
$ws.Range("H5").PasteSpecial(-4122) | Out-Null

# 4) Fill in the actual figures.
$ws.Range("H3").Value = 2020

$ws.Range("G4").Value = 146
$ws.Range("H4").Value = 158

$ws.Range("G5").Value = 127
$ws.Range("H5").Value = 397

$excel.CutCopyMode = $false
